# Actualización automática de tasas-transfi.xlsx
$wb = $excel.ActiveWorkbook

# --- Hoja1: update the "Conversión del día" text block ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")

$newText = "Conversión del día 💰" + [char]10 + `
    "✅ Dólar paralelo: 68" + [char]10 + `
    [char]10 + `
    "Binance" + [char]10 + `
    "✅ 1000 Bs = 3.32 = 12730.9 pesos" + [char]10 + `
    "✅ 12730.9 pesos = 3.3 = 957.96 Bs" + [char]10 + `
    [char]10 + `
    "Promedio competencia" + [char]10 + `
    "✅ Tasa pesos: 20" + [char]10 + `
    "✅ Tasa Bs: 20" + [char]10 + `
    "✅ % Ganancia: 20%"

$wsHoja1.Range("A1").Value = $newText

# --- tasas: update rate cells ---
$wsTasas = $wb.Worksheets.Item("tasas")

$wsTasas.Range("N10").Value = 301
$wsTasas.Range("O10").Value = 3832
$wsTasas.Range("N12").Value = 3854
$wsTasas.Range("O12").Value = 290
